$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new value looks like a plain number but must
# stay stored as text (matches the source feed formatting, e.g. "307.36").
# Force text format before assigning so Excel does not coerce them to numbers.
$textCells = @(
    "D5",
    "D6",
    "D7",
    "D9",
    "D10",
    "D11",
    "D13",
    "D15",
    "D17",
    "D19",
    "D21",
    "D22",
    "D23",
    "D25",
    "D28",
    "D30",
    "D31",
    "D33",
    "D36",
    "D37",
    "D38",
    "D39",
    "D44",
    "D45",
    "D46",
    "D47",
    "D48",
    "D50"
)
foreach ($cell in $textCells) {
    $ws.Range($cell).NumberFormat = "@"
}

# Updated coin prices / 1h volume changes (and two ranking swaps: rows
# 38/39 Kaspa<->Stellar, rows 44/45 EnergySwap<->ApeXProtocol).
$ws.Range("D2").Value = '42.082.65'
$ws.Range("E2").Value = '  -0.73%  '
$ws.Range("D3").Value = '2.254.24'
$ws.Range("E3").Value = '  -0.96%  '
$ws.Range("D5").Value = '307.36'
$ws.Range("E5").Value = '  -0.01%  '
$ws.Range("D6").Value = '96.78'
$ws.Range("E6").Value = '  -1.33%  '
$ws.Range("D7").Value = '0.524'
$ws.Range("E7").Value = '  -0.97%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("D9").Value = '0.487'
$ws.Range("E9").Value = '  -1.63%  '
$ws.Range("D10").Value = '34.63'
$ws.Range("E10").Value = '  -3.36%  '
$ws.Range("D11").Value = '0.0827'
$ws.Range("E11").Value = '  +3.75%  '
$ws.Range("E12").Value = '  +0.47%  '
$ws.Range("D13").Value = '6.79'
$ws.Range("E13").Value = '  +1.47%  '
$ws.Range("D14").Value = '2.604.57'
$ws.Range("E14").Value = '  -0.95%  '
$ws.Range("D15").Value = '14.55'
$ws.Range("E15").Value = '  +0.93%  '
$ws.Range("D16").Value = '2.256.11'
$ws.Range("E16").Value = '  -0.80%  '
$ws.Range("D17").Value = '0.782'
$ws.Range("E17").Value = '  -1.53%  '
$ws.Range("D18").Value = '41.958.65'
$ws.Range("E18").Value = '  -0.82%  '
$ws.Range("D19").Value = '12.24'
$ws.Range("E19").Value = '  -2.29%  '
$ws.Range("D20").Value = '0.0₃0906'
$ws.Range("E20").Value = '  -0.48%  '
$ws.Range("D21").Value = '5.93'
$ws.Range("E21").Value = '  -0.86%  '
$ws.Range("D22").Value = '67.32'
$ws.Range("E22").Value = '  -0.43%  '
$ws.Range("D23").Value = '235.17'
$ws.Range("E23").Value = '  -2.39%  '
$ws.Range("E24").Value = '  -1.09%  '
$ws.Range("D25").Value = '1.95'
$ws.Range("E25").Value = '  +0.50%  '
$ws.Range("E26").Value = '  +0.00%  '
$ws.Range("E27").Value = '  -1.61%  '
$ws.Range("D28").Value = '36.87'
$ws.Range("E28").Value = '  -2.21%  '
$ws.Range("E29").Value = '  +1.01%  '
$ws.Range("D30").Value = '9.51'
$ws.Range("E30").Value = '  +0.24%  '
$ws.Range("D31").Value = '164.96'
$ws.Range("E31").Value = '  +3.70%  '
$ws.Range("E32").Value = '  -0.04%  '
$ws.Range("D33").Value = '5.19'
$ws.Range("E33").Value = '  -0.60%  '
$ws.Range("E34").Value = '  -1.68%  '
$ws.Range("E35").Value = '  +2.90%  '
$ws.Range("D36").Value = '0.0724'
$ws.Range("E36").Value = '  -2.43%  '
$ws.Range("D37").Value = '2.37'
$ws.Range("E37").Value = '  -0.46%  '
$ws.Range("B38").Value = 'Stellar'
$ws.Range("C38").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D38").Value = '0.115'
$ws.Range("E38").Value = '  -0.07%  '
$ws.Range("B39").Value = 'Kaspa'
$ws.Range("C39").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D39").Value = '0.103'
$ws.Range("E39").Value = '  -2.20%  '
$ws.Range("E40").Value = '  -2.25%  '
$ws.Range("E41").Value = '  +0.61%  '
$ws.Range("D42").Value = '1.936.99'
$ws.Range("E42").Value = '  -3.05%  '
$ws.Range("E43").Value = '  -1.43%  '
$ws.Range("B44").Value = 'ApeXProtocol'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D44").Value = '2.20'
$ws.Range("E44").Value = '  -8.62%  '
$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").Value = '18.52'
$ws.Range("E45").Value = '  -2.95%  '
$ws.Range("D46").Value = '2.92'
$ws.Range("E46").Value = '  -2.54%  '
$ws.Range("D47").Value = '9.67'
$ws.Range("E47").Value = '  -3.03%  '
$ws.Range("D48").Value = '53.81'
$ws.Range("E48").Value = '  +1.41%  '
$ws.Range("D49").Value = '2.477.79'
$ws.Range("E49").Value = '  -0.81%  '
$ws.Range("D50").Value = '71.38'
$ws.Range("E50").Value = '  -1.11%  '
$ws.Range("E51").Value = '  -0.57%  '
